$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 54.53585066666667
$ws.Range("H2").Value = 163.607552
$ws.Range("I2").Value = 0.3031388658437607
$ws.Range("J2").Value = 0.3031388658437607
$ws.Range("M2").Value = 13.713764
$ws.Range("N2").Value = 41.141292
$ws.Range("O2").Value = 0.0901423721847377
$ws.Range("P2").Value = 0.0901423721847377
$ws.Range("Q2").Value = 747.8917855819093
$ws.Range("R2").Value = 6731.026070237184
$ws.Range("S2").Value = 0.02732565646854754
$ws.Range("T2").Value = 0.02732565646854755
# Row 3
$ws.Range("G3").Value = 54.53585066666667
$ws.Range("H3").Value = 163.607552
$ws.Range("I3").Value = 0.3031388658437607
$ws.Range("J3").Value = 0.3031388658437607
$ws.Range("N3").Value = 84.55600199999999
$ws.Range("O3").Value = 0.1852659027513629
$ws.Range("P3").Value = 0.1852659027513629
$ws.Range("Q3").Value = 1537.111166014123
$ws.Range("R3").Value = 13834.0004941271
$ws.Range("S3").Value = 0.05616129563956861
$ws.Range("T3").Value = 0.05616129563956861
# Row 4
$ws.Range("G4").Value = 54.53585066666667
$ws.Range("H4").Value = 163.607552
$ws.Range("I4").Value = 0.3031388658437607
$ws.Range("J4").Value = 0.3031388658437607
$ws.Range("M4").Value = 21.07704566666666
$ws.Range("N4").Value = 63.23113699999999
$ws.Range("O4").Value = 0.1385421898057586
$ws.Range("P4").Value = 0.1385421898057586
$ws.Range("Q4").Value = 1149.454614971847
$ws.Range("R4").Value = 10345.09153474662
$ws.Range("S4").Value = 0.04199752228922869
$ws.Range("T4").Value = 0.04199752228922869
# Row 5
$ws.Range("G5").Value = 54.53585066666667
$ws.Range("H5").Value = 163.607552
$ws.Range("I5").Value = 0.3031388658437607
$ws.Range("J5").Value = 0.3031388658437607
$ws.Range("M5").Value = 89.15834833333334
$ws.Range("N5").Value = 267.475045
$ws.Range("O5").Value = 0.5860495352581409
$ws.Range("P5").Value = 0.5860495352581407
$ws.Range("Q5").Value = 4862.326370393316
$ws.Range("R5").Value = 43760.93733353984
$ws.Range("S5").Value = 0.1776543914464158
$ws.Range("T5").Value = 0.1776543914464158
# Row 6
$ws.Range("I6").Value = 0.1026363515063155
$ws.Range("J6").Value = 0.1026363515063155
$ws.Range("M6").Value = 13.713764
$ws.Range("N6").Value = 41.141292
$ws.Range("O6").Value = 0.0901423721847377
$ws.Range("P6").Value = 0.0901423721847377
$ws.Range("Q6").Value = 253.2201998579546
$ws.Range("R6").Value = 2278.981798721592
$ws.Range("S6").Value = 0.009251884197165855
$ws.Range("T6").Value = 0.009251884197165855
# Row 7
$ws.Range("I7").Value = 0.1026363515063155
$ws.Range("J7").Value = 0.1026363515063155
$ws.Range("N7").Value = 84.55600199999999
$ws.Range("O7").Value = 0.1852659027513629
$ws.Range("P7").Value = 0.1852659027513629
$ws.Range("R7").Value = 4683.897373244051
$ws.Range("S7").Value = 0.01901501631692374
$ws.Range("T7").Value = 0.01901501631692375
# Row 8
$ws.Range("I8").Value = 0.1026363515063155
$ws.Range("J8").Value = 0.1026363515063155
$ws.Range("M8").Value = 21.07704566666666
$ws.Range("N8").Value = 63.23113699999999
$ws.Range("O8").Value = 0.1385421898057586
$ws.Range("P8").Value = 0.1385421898057586
$ws.Range("Q8").Value = 389.1808052208401
$ws.Range("R8").Value = 3502.627246987561
$ws.Range("S8").Value = 0.01421946489135852
$ws.Range("T8").Value = 0.01421946489135852
# Row 9
$ws.Range("I9").Value = 0.1026363515063155
$ws.Range("J9").Value = 0.1026363515063155
$ws.Range("M9").Value = 89.15834833333334
$ws.Range("N9").Value = 267.475045
$ws.Range("O9").Value = 0.5860495352581409
$ws.Range("P9").Value = 0.5860495352581407
$ws.Range("Q9").Value = 1646.279955231241
$ws.Range("R9").Value = 14816.51959708117
$ws.Range("S9").Value = 0.06014998610086738
$ws.Range("T9").Value = 0.06014998610086737
# Row 10
$ws.Range("G10").Value = 12.55635966666667
$ws.Range("H10").Value = 37.669079
$ws.Range("I10").Value = 0.06979483370938171
$ws.Range("J10").Value = 0.06979483370938172
$ws.Range("M10").Value = 13.713764
$ws.Range("N10").Value = 41.141292
$ws.Range("O10").Value = 0.0901423721847377
$ws.Range("P10").Value = 0.0901423721847377
$ws.Range("Q10").Value = 172.1949531677853
$ws.Range("R10").Value = 1549.754578510068
$ws.Range("S10").Value = 0.006291471876802963
$ws.Range("T10").Value = 0.006291471876802963
# Row 11
$ws.Range("G11").Value = 12.55635966666667
$ws.Range("H11").Value = 37.669079
$ws.Range("I11").Value = 0.06979483370938171
$ws.Range("J11").Value = 0.06979483370938172
$ws.Range("N11").Value = 84.55600199999999
$ws.Range("O11").Value = 0.1852659027513629
$ws.Range("P11").Value = 0.1852659027513629
$ws.Range("Q11").Value = 353.9051910291286
$ws.Range("R11").Value = 3185.146719262158
$ws.Range("S11").Value = 0.01293060287454986
$ws.Range("T11").Value = 0.01293060287454986
# Row 12
$ws.Range("G12").Value = 12.55635966666667
$ws.Range("H12").Value = 37.669079
$ws.Range("I12").Value = 0.06979483370938171
$ws.Range("J12").Value = 0.06979483370938172
$ws.Range("M12").Value = 21.07704566666666
$ws.Range("N12").Value = 63.23113699999999
$ws.Range("O12").Value = 0.1385421898057586
$ws.Range("P12").Value = 0.1385421898057586
$ws.Range("Q12").Value = 264.6509661014247
$ws.Range("R12").Value = 2381.858694912823
$ws.Range("S12").Value = 0.009669529099226519
$ws.Range("T12").Value = 0.009669529099226521
# Row 13
$ws.Range("G13").Value = 12.55635966666667
$ws.Range("H13").Value = 37.669079
$ws.Range("I13").Value = 0.06979483370938171
$ws.Range("J13").Value = 0.06979483370938172
$ws.Range("M13").Value = 89.15834833333334
$ws.Range("N13").Value = 267.475045
$ws.Range("O13").Value = 0.5860495352581409
$ws.Range("P13").Value = 0.5860495352581407
$ws.Range("Q13").Value = 1119.504288959284
$ws.Range("R13").Value = 10075.53860063356
$ws.Range("S13").Value = 0.04090322985880237
$ws.Range("T13").Value = 0.04090322985880238
# Row 14
$ws.Range("G14").Value = 94.34696966666667
$ws.Range("H14").Value = 283.040909
$ws.Range("I14").Value = 0.524429948940542
$ws.Range("J14").Value = 0.5244299489405421
$ws.Range("M14").Value = 13.713764
$ws.Range("N14").Value = 41.141292
$ws.Range("O14").Value = 0.0901423721847377
$ws.Range("P14").Value = 0.0901423721847377
$ws.Range("Q14").Value = 1293.852076123825
$ws.Range("R14").Value = 11644.66868511443
$ws.Range("S14").Value = 0.04727335964222132
$ws.Range("T14").Value = 0.04727335964222133
# Row 15
$ws.Range("G15").Value = 94.34696966666667
$ws.Range("H15").Value = 283.040909
$ws.Range("I15").Value = 0.524429948940542
$ws.Range("J15").Value = 0.5244299489405421
$ws.Range("N15").Value = 84.55600199999999
$ws.Range("O15").Value = 0.1852659027513629
$ws.Range("P15").Value = 0.1852659027513629
$ws.Range("Q15").Value = 2659.200851942868
$ws.Range("R15").Value = 23932.80766748582
$ws.Range("S15").Value = 0.09715898792032066
$ws.Range("T15").Value = 0.09715898792032068
# Row 16
$ws.Range("G16").Value = 94.34696966666667
$ws.Range("H16").Value = 283.040909
$ws.Range("I16").Value = 0.524429948940542
$ws.Range("J16").Value = 0.5244299489405421
$ws.Range("M16").Value = 21.07704566666666
$ws.Range("N16").Value = 63.23113699999999
$ws.Range("O16").Value = 0.1385421898057586
$ws.Range("P16").Value = 0.1385421898057586
$ws.Range("Q16").Value = 1988.555388175948
$ws.Range("R16").Value = 17896.99849358353
$ws.Range("S16").Value = 0.07265567352594486
$ws.Range("T16").Value = 0.07265567352594487
# Row 17
$ws.Range("G17").Value = 94.34696966666667
$ws.Range("H17").Value = 283.040909
$ws.Range("I17").Value = 0.524429948940542
$ws.Range("J17").Value = 0.5244299489405421
$ws.Range("M17").Value = 89.15834833333334
$ws.Range("N17").Value = 267.475045
$ws.Range("O17").Value = 0.5860495352581409
$ws.Range("P17").Value = 0.5860495352581407
$ws.Range("Q17").Value = 8411.819985735101
$ws.Range("R17").Value = 75706.37987161591
$ws.Range("S17").Value = 0.3073419278520552
$ws.Range("T17").Value = 0.3073419278520552

Write-Host "done"